$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 92
$ws1.Range("F4").Value = 112
$ws1.Range("F5").Value = 45
$ws1.Range("F6").Value = 68
$ws1.Range("F9").Value = 1112
$ws1.Range("F10").Value = 277
$ws1.Range("F11").Value = 4
$ws1.Range("F12").Value = 10348
$ws1.Range("F13").Value = 3
$ws1.Range("F15").Value = 275
$ws1.Range("F16").Value = 19
$ws1.Range("F17").Value = 656
$ws1.Range("F18").Value = 11900
$ws1.Range("F19").Value = 12286
$ws1.Range("F20").Value = 28
$ws1.Range("F21").Value = 109
$ws1.Range("F24").Value = 71
$ws1.Range("F25").Value = 36

# Sheet "全部类型" (sheet index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 92
$ws4.Range("F4").Value = 112
$ws4.Range("F5").Value = 45
$ws4.Range("F6").Value = 68
$ws4.Range("F10").Value = 1112
$ws4.Range("F11").Value = 278
$ws4.Range("F12").Value = 4
$ws4.Range("F13").Value = 10348
$ws4.Range("F14").Value = 3
$ws4.Range("F16").Value = 275
$ws4.Range("F17").Value = 19
$ws4.Range("F18").Value = 656
$ws4.Range("F19").Value = 11900
$ws4.Range("F20").Value = 12286
$ws4.Range("F21").Value = 28
$ws4.Range("F22").Value = 109
$ws4.Range("F25").Value = 71
$ws4.Range("F26").Value = 36
